# Agate - indicateurs statistiques.xlsx
# "Amelioration du calcul statistique et de l'affichage"
#
# 1. "type indicateur" sheet: fix the "superficie" label (m2 -> km2) and add
#    two new rows describing the zoning fields (idZonage / idZonage.name).
# 2. Add a brand new "Zone predefinie" sheet listing the predefined zoning
#    levels (Departement, Commune, QPV, Zonage utilisateur).
# 3. "categorie" sheet: extend the AutoFilter from A1:C27 to A1:D27, which
#    swaps which _xlnm._FilterDatabase defined name is the hidden one.

$wb = $excel.ActiveWorkbook

# --- 1. "type indicateur" sheet -------------------------------------------
$wsType = $wb.Worksheets.Item("type indicateur")

# Fix the unit in the existing "superficie" label.
$wsType.Range("B30").Value = "Superficie (en km²)"

# New rows documenting the zoning columns used elsewhere in the workbook.
$wsType.Range("A32").Value = "idZonage"
$wsType.Range("B32").Value = "Identifiant de la zone"
$wsType.Range("A33").Value = "idZonage.name"
$wsType.Range("B33").Value = "Zone"

# --- 2. New "Zone predefinie" sheet ----------------------------------------
$sheetCount = $wb.Worksheets.Count
$wsZone = $wb.Worksheets.Add($null, $wb.Worksheets.Item($sheetCount))
$wsZone.Name = "Zone predefinie"

$wsZone.Range("A1").Value = "idPredefine"
$wsZone.Range("B1").Value = "labelPredefine"
$wsZone.Range("A2").Value = 1
$wsZone.Range("B2").Value = "Département"
$wsZone.Range("A3").Value = 2
$wsZone.Range("B3").Value = "Commune"
$wsZone.Range("A4").Value = 3
$wsZone.Range("B4").Value = "QPV"
$wsZone.Range("A5").Value = 4
$wsZone.Range("B5").Value = "Zonage utilisateur"

# --- 3. "categorie" sheet: widen the AutoFilter to include column D --------
$wsCat = $wb.Worksheets.Item("categorie")
$wsCat.AutoFilterMode = $false
$null = $wsCat.Range("A1:D27").AutoFilter()

# Re-apply creates/updates the _xlnm._FilterDatabase defined names but
# leaves the old (C27) range marked visible and the new (D27) one hidden;
# the canonical file keeps it the other way around, so fix that up.
$n1 = $wb.Names.Item(1)
$n3 = $wb.Names.Item(3)
$n1.RefersTo = "=categorie!`$A`$1:`$D`$27"
$n1.Visible = $false
$n3.RefersTo = "=categorie!`$A`$1:`$C`$27"
$n3.Visible = $true

# --- Final selection / active sheet ----------------------------------------
$null = $wsZone.Range("B5").Select()
$null = $wsType.Range("A33").Select()
